# Append a new results row (row 30) to the "Results" sheet, duplicating the
# previous row (29) of M3C2 statistics but with an updated run Timestamp.
# This mirrors a new pipeline run ("only_stats" config) that produced an
# identical statistics/inlier-outlier summary with a later timestamp.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

$row = 30

$ws.Cells.Item($row, 1).Value = "2025-08-26 14:16:15"  # A30
$ws.Cells.Item($row, 2).Value = "data\0342-0349"  # B30
$ws.Cells.Item($row, 3).Value = "ref"  # C30
$ws.Cells.Item($row, 4).Value = 709128  # D30
$ws.Cells.Item($row, 5).Value = 0.1245588149878983  # E30
$ws.Cells.Item($row, 6).Value = 0.2491176299757966  # F30
$ws.Cells.Item($row, 7).Value = 95  # G30
$ws.Cells.Item($row, 8).Value = 0.00013396735145136  # H30
$ws.Cells.Item($row, 9).Value = 0.9998660326485487  # I30
$ws.Cells.Item($row, 10).Value = 709033  # J30
$ws.Cells.Item($row, 11).Value = -151.9928180000001  # K30
$ws.Cells.Item($row, 12).Value = 151.729350141048  # L30
$ws.Cells.Item($row, 13).Value = 690813  # M30
$ws.Cells.Item($row, 14).Value = -1125.949142  # N30
$ws.Cells.Item($row, 15).Value = 62.170601304942  # O30
$ws.Cells.Item($row, 16).Value = -0.120008  # P30
$ws.Cells.Item($row, 17).Value = 0.134913  # Q30
$ws.Cells.Item($row, 18).Value = -0.0002143663524828888  # R30
$ws.Cells.Item($row, 19).Value = -0.002299  # S30
$ws.Cells.Item($row, 20).Value = 0.01462855990518499  # T30
$ws.Cells.Item($row, 21).Value = 0.01462698916272617  # U30
$ws.Cells.Item($row, 22).Value = 0.008219008269008636  # V30
$ws.Cells.Item($row, 23).Value = 0.006827373  # W30
$ws.Cells.Item($row, 24).Value = -0.043862  # X30
$ws.Cells.Item($row, 25).Value = 0.043885  # Y30
$ws.Cells.Item($row, 26).Value = -0.00162988991521584  # Z30
$ws.Cells.Item($row, 27).Value = -0.002469  # AA30
$ws.Cells.Item($row, 28).Value = 0.009486637003424522  # AB30
$ws.Cells.Item($row, 29).Value = 0.009345573310344375  # AC30
$ws.Cells.Item($row, 30).Value = 0.006657630840762984  # AD30
$ws.Cells.Item($row, 31).Value = 0.0065990526  # AE30
$ws.Cells.Item($row, 32).Value = 690813  # AF30
$ws.Cells.Item($row, 33).Value = 246292  # AG30
$ws.Cells.Item($row, 34).Value = 444509  # AH30
$ws.Cells.Item($row, 35).Value = 15929  # AI30
$ws.Cells.Item($row, 36).Value = 2291  # AJ30
$ws.Cells.Item($row, 37).Value = 18220  # AK30
$ws.Cells.Item($row, 38).Value = 0.05345534160263447  # AL30
$ws.Cells.Item($row, 39).Value = 0.04536447157847037  # AM30
$ws.Cells.Item($row, 40).Value = -0.013568  # AN30
$ws.Cells.Item($row, 41).Value = -0.006508  # AO30
$ws.Cells.Item($row, 42).Value = 0.002725  # AP30
$ws.Cells.Item($row, 43).Value = 0.022282  # AQ30
$ws.Cells.Item($row, 44).Value = 0.009233  # AR30
$ws.Cells.Item($row, 45).Value = -0.013419  # AS30
$ws.Cells.Item($row, 46).Value = -0.006579  # AT30
$ws.Cells.Item($row, 47).Value = 0.002358  # AU30
$ws.Cells.Item($row, 48).Value = 0.013891  # AV30
$ws.Cells.Item($row, 49).Value = 0.008937  # AW30
$ws.Cells.Item($row, 50).Value = -0.0002143663524828888  # AX30
$ws.Cells.Item($row, 51).Value = 0.01462698916272617  # AY30
$ws.Cells.Item($row, 52).Value = 32673381941410.1  # AZ30
$ws.Cells.Item($row, 53).Value = 1.010085389383569  # BA30
$ws.Cells.Item($row, 54).Value = 0.3498892426805521  # BB30
$ws.Cells.Item($row, 55).Value = -0.09252700291488647  # BC30
$ws.Cells.Item($row, 56).Value = -0.08886902402080391  # BD30
$ws.Cells.Item($row, 57).Value = 1.970131977073834  # BE30
$ws.Cells.Item($row, 58).Value = 11637763.60543192  # BF30
$ws.Cells.Item($row, 59).Value = 2.908214737699866  # BG30
$ws.Cells.Item($row, 60).Value = 16.58429873827257  # BH30
$ws.Cells.Item($row, 61).Value = "data\0342-0349\python_ref_m3c2_distances.txt"  # BI30
$ws.Cells.Item($row, 62).Value = "data\0342-0349\python_ref_m3c2_params.txt"  # BJ30
$ws.Cells.Item($row, 63).Value = 3  # BK30
